# Update cryptocurrency price/volume data to reflect the latest GitHub Actions refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin), C (Link), D (Price) and E (Volume(1h)) are plain-text cells in this
# sheet (t="inlineStr"); force text format on D/E so Excel keeps the numeric-looking
# strings ("310.57", "0.73%", ...) as text instead of auto-converting them to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.57"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.73%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.55"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.45%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.112"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.07%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07804"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.97%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.908"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.81%"

$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.218"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.14%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.914"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-8.24%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9263"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.20%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1194"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.50%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1903"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.29%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09419"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "6.51%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03438"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.47%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09618"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.38%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001364"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.91%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005902"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.54%"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.86%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.401"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.47%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3425"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.33%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.264"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.18%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.23%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2591"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.83%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.02103"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "180.20%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04349"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.60%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001199"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.85%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004257"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.56%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-63.80%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02076"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-9.55%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05080"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.30%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007641"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.13%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009131"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-7.62%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1348"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.42%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002071"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.02%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.30%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006712"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.45%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.13%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002915"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.91%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001201"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.13%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.13%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.13%"
